# Daily rollover update:
#   For every data row (2..99) on the active sheet:
#     - Column D = total cycle length (days)
#     - Column E = days remaining in the current cycle
#     - Column F = start date (yyyyMMdd) of the current cycle
#   Each day, the "remaining" counter ticks down by 1. When it would drop
#   to 0 (i.e. it was 1), the cycle rolls over: remaining resets to the
#   full cycle length (D) and the start date advances by D days.
#   Rows whose F value isn't a well-formed yyyyMMdd date are left
#   untouched (mirrors source data that failed to parse upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Range("A1").End(4).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string][int]$fVal

    if ($fStr.Length -ne 8) {
        # Not a valid yyyyMMdd date (e.g. a data-entry typo) - skip, as
        # the source update did for this row.
        continue
    }

    $parsedOk = $true
    try {
        $startDate = [DateTime]::ParseExact($fStr, "yyyyMMdd", $null)
    } catch {
        $parsedOk = $false
    }

    if (-not $parsedOk) {
        continue
    }

    $remaining = [int]$eVal
    $total = [int]$dVal

    if ($remaining -le 1) {
        $newRemaining = $total
        $newStart = $startDate.AddDays($total)
        $newFVal = [int]$newStart.ToString("yyyyMMdd")
        $eCell.Value = $newRemaining
        $fCell.Value = $newFVal
    } else {
        $newRemaining = $remaining - 1
        $eCell.Value = $newRemaining
    }
}
